# Update each two-digit-divided-by-one-digit answer cell with new values.
# Replacements are applied in document order, replacing only the first
# remaining occurrence each time (wdReplaceOne) so that duplicate original
# texts (e.g. "72÷8=9, 0" appearing twice) are each mapped to their own
# distinct new value, and no replacement output is re-matched by a later
# find (since the cursor only advances forward through the document).

$d = $word.ActiveDocument

$pairs = @(
    @("60÷2=30, 0", "45÷5=9, 0"),
    @("64÷9=7, 1", "48÷7=6, 6"),
    @("72÷8=9, 0", "86÷8=10, 6"),
    @("72÷8=9, 0", "38÷9=4, 2"),
    @("57÷6=9, 3", "62÷9=6, 8"),
    @("30÷4=7, 2", "71÷2=35, 1"),
    @("97÷6=16, 1", "83÷7=11, 6"),
    @("72÷7=10, 2", "11÷2=5, 1"),
    @("34÷9=3, 7", "35÷4=8, 3"),
    @("91÷2=45, 1", "92÷9=10, 2"),
    @("67÷9=7, 4", "10÷8=1, 2"),
    @("49÷2=24, 1", "72÷2=36, 0"),
    @("53÷7=7, 4", "32÷3=10, 2"),
    @("92÷2=46, 0", "99÷7=14, 1"),
    @("59÷7=8, 3", "67÷5=13, 2"),
    @("35÷2=17, 1", "73÷3=24, 1"),
    @("78÷2=39, 0", "68÷6=11, 2"),
    @("34÷5=6, 4", "57÷2=28, 1"),
    @("21÷2=10, 1", "38÷5=7, 3"),
    @("84÷7=12, 0", "42÷4=10, 2"),
    @("87÷5=17, 2", "19÷9=2, 1"),
    @("30÷8=3, 6", "37÷5=7, 2"),
    @("47÷9=5, 2", "59÷7=8, 3"),
    @("19÷5=3, 4", "15÷8=1, 7"),
    @("59÷4=14, 3", "28÷8=3, 4")
)

$rng = $d.Content

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
    # Collapse to the end of the replaced text so the next search starts
    # just after it, preserving document order and not re-matching it.
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}
